$wb = $excel.ActiveWorkbook

# Update selection on "中信兄弟" sheet (view-state only change)
$sheet2 = $wb.Worksheets.Item("中信兄弟")
$sheet2.Range("B22").Select()

# "Dragons" sheet loses tabSelected as it is no longer the active tab
$sheet4 = $wb.Worksheets.Item("Dragons")
$sheet4.Activate()

# The sheet that used to be "統一7-ELEVEn獅" keeps its data but gets a new
# selection and is no longer the active tab once the new sheet is inserted
$oldLast = $wb.Worksheets.Item("統一7-ELEVEn獅")
$oldLast.Range("N21").Select()

# Insert a new worksheet "台鋼雄鷹" right before "統一7-ELEVEn獅"
$newSheet = $wb.Worksheets.Add($oldLast)
$newSheet.Name = "台鋼雄鷹"

$headers = @("年度","出賽數","勝","敗","和","勝率","主場勝","主場敗","主場和","客場勝","客場敗","客場和")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$newSheet.Range("A2").Value = "2023(下)"
$newSheet.Range("A3").Value = "2023(上)"

$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

$newSheet.Range("A3").Select()
